$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix sign of formula in C8 (+C4/D4 -> -C4/D4)
# ---------------------------------------------------------------------------
$ws.Range("C8").Formula = "=-C4/D4"

# ---------------------------------------------------------------------------
# 2) Apply the orange fill style (same style already used on B2/D4/etc, the
#    "s=2" cellXf) to the alfa-summary labels B10, B11, B12
# ---------------------------------------------------------------------------
$fillColor = $ws.Range("B2").Interior.Color
$ws.Range("B10").Interior.Color = $fillColor
$ws.Range("B11").Interior.Color = $fillColor
$ws.Range("B12").Interior.Color = $fillColor

# ---------------------------------------------------------------------------
# 3) Extend the iteration-index header row (row 17) from column H out to K
# ---------------------------------------------------------------------------
$ws.Range("I17").Interior.Color = $fillColor
$ws.Range("I17").Value = 7
$ws.Range("J17").Interior.Color = $fillColor
$ws.Range("J17").Value = 8
$ws.Range("K17").Interior.Color = $fillColor
$ws.Range("K17").Value = 9

# ---------------------------------------------------------------------------
# 4) Jacobi/Gauss-Seidel iteration rows 18-20: change the row 19/20 formulas
#    to reference the just-updated column (Gauss-Seidel style) instead of the
#    previous column, and extend all three rows from column H out to K.
# ---------------------------------------------------------------------------
$ws.Range("D18:K18").Formula = "=+`$C`$6*C19+`$D`$6*C20+`$F`$6"

$ws.Range("C19").Formula = "=+`$B`$7*C18+`$D`$7*B20+`$F`$7"
$ws.Range("D19:K19").Formula = "=+`$B`$7*D18+`$D`$7*C20+`$F`$7"

$ws.Range("C20").Formula = "=+`$B`$8*C18+`$C`$8*C19+`$F`$8"
$ws.Range("D20:K20").Formula = "=+`$B`$8*D18+`$C`$8*D19+`$F`$8"

# ---------------------------------------------------------------------------
# 5) New error rows (e1, e2, e3) computing the |difference| between
#    successive iterations for x1, x2, x3 respectively.
# ---------------------------------------------------------------------------
$ws.Range("A23").Interior.Color = $fillColor
$ws.Range("A23").Value = "e1"
$ws.Range("B23").Formula = "=+ABS(C18)-ABS(B18)"
$ws.Range("C23:I23").Formula = "=+ABS(D18)-ABS(C18)"

$ws.Range("A24").Interior.Color = $fillColor
$ws.Range("A24").Value = "e2"
$ws.Range("B24:I25").Formula = "=+ABS(C19)-ABS(B19)"

$ws.Range("A25").Interior.Color = $fillColor
$ws.Range("A25").Value = "e3"
$ws.Range("I25").Interior.Color = $fillColor

# ---------------------------------------------------------------------------
# 6) Row with the max error of each iteration column.
# ---------------------------------------------------------------------------
$ws.Range("B27").Formula = "=+MAX(B23:B25)"
$ws.Range("C27:I27").Formula = "=+MAX(C23:C25)"

# ---------------------------------------------------------------------------
# 7) Window view: zoom to 70% and move the selection to K28 (no frozen
#    top-left scroll position anymore).
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Zoom = 70
$ws.Range("K28").Select() | Out-Null
